$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing row 57 (and all below) down to 58.
$ws.Rows(57).Insert()

# Populate the newly inserted row 57 with the new record's data.
$ws.Cells.Item(57, 1).Value = 7
$ws.Cells.Item(57, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(57, 3).Value = "Ñuble"
$ws.Cells.Item(57, 4).Value = 45086
$ws.Cells.Item(57, 5).Value = 16
$ws.Cells.Item(57, 6).Value = 100112001
$ws.Cells.Item(57, 7).Value = "Berenjena"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 100
$ws.Cells.Item(57, 11).Value = 9000
$ws.Cells.Item(57, 12).Value = 10000
$ws.Cells.Item(57, 13).Value = 9500
$ws.Cells.Item(57, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(57, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(57, 16).Value = 158
$ws.Cells.Item(57, 17).Value = 60
$ws.Cells.Item(57, 18).Value = "Hortaliza"
